# Update the build timestamp embedded in the "version" strings across the
# workbook, from "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST".

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# "About" sheet: version banner (A2) and recommended citation (A6).
$targets = @(
    @{ Sheet = $aboutSheet; Cell = "A2" },
    @{ Sheet = $aboutSheet; Cell = "A6" },
    @{ Sheet = $dataSheet; Cell = "S2" },
    @{ Sheet = $dataSheet; Cell = "S3" },
    @{ Sheet = $dataSheet; Cell = "S4" },
    @{ Sheet = $dataSheet; Cell = "S5" },
    @{ Sheet = $dataSheet; Cell = "S6" },
    @{ Sheet = $dataSheet; Cell = "S7" }
)

foreach ($t in $targets) {
    $cell = $t.Sheet.Range($t.Cell)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [string] -and $val.Contains($oldStamp)) {
        $cell.Value2 = $val.Replace($oldStamp, $newStamp)
    }
}
